$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: copy formatting (fill/font/alignment) from a reference cell
#     that already carries the desired status style, then set the new text.
function Set-StatusCell($destRef, $srcRef, $text) {
    $ws.Range($srcRef).Copy() | Out-Null
    $ws.Range($destRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($destRef).Value = $text
}

# Reference cells that already use the canonical status formatting:
#   TODO         -> F8
#   DONE         -> F4
#   IN PROGRESS  -> G20

# F5: TODO -> IN PROGRESS
Set-StatusCell "F5" "G20" "IN PROGRESS"

# F6: DONE -> IN PROGRESS
Set-StatusCell "F6" "G20" "IN PROGRESS"

# F22: DONE -> TODO
Set-StatusCell "F22" "F8" "TODO"

# F37: DONE -> IN PROGRESS
Set-StatusCell "F37" "G20" "IN PROGRESS"

# F38: DONE -> IN PROGRESS
Set-StatusCell "F38" "G20" "IN PROGRESS"

# F49: IN PROGRESS: EN -> TODO
Set-StatusCell "F49" "F8" "TODO"

# F62: IN PROGRESS: EN -> TODO
Set-StatusCell "F62" "F8" "TODO"

# F63: IN PROGRESS: EN -> TODO
Set-StatusCell "F63" "F8" "TODO"

# H35: remove stale note entirely (Clear, not just ClearContents, so the
# cell element itself drops out of the sheet rather than lingering as an
# empty-but-styled cell)
$ws.Range("H35").Clear() | Out-Null

# H54: update note text (format unchanged)
$ws.Range("H54").Value = "Require an integer, followed by anything at all. We _don't_ validate semanticVersion, which should be a separate check in the same rule, and which in fact could have been in the schema!"

# Reset the sheet view: drop the custom zoom/top-left scroll position and
# move the active selection to H55.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H55").Select() | Out-Null
